$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.733288645744324
$ws.Range("B1").Value = 4.094292163848877
$ws.Range("C1").Value = 3.195954084396362
$ws.Range("D1").Value = 1.599943518638611
$ws.Range("E1").Value = 0.7425632476806641
